$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new "Wins" / "Losses" / "Ties" columns (AD:AF), matching the
# bold/bordered header style already used by the existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Season record for every player row: 53 wins, 62 losses, 0 ties.
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 30).Value = 53
    $ws.Cells.Item($row, 31).Value = 62
    $ws.Cells.Item($row, 32).Value = 0
}
